# Updated symbol list on Sun Dec 18 13:19:08 UTC 2022 with GitHub Actions
#
# Refreshes the "cryptos" price sheet: most Price (column D) cells get a
# new quote, and the coins that moved in/out of the top-50 ranking have
# their whole row (Coin / Link / Price / Volume(1h)) rewritten in place.
# Price values must stay TEXT (the sheet stores "247.49" etc. as strings,
# not numbers), so every Price write goes through Set-TextValue, which
# forces a Text format before the assignment and then restores the
# original "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# --- Price-only refreshes ---
Set-TextValue $ws.Cells.Item(2, 4) "246.87"
Set-TextValue $ws.Cells.Item(3, 4) "22.38"
Set-TextValue $ws.Cells.Item(5, 4) "0.05633"
Set-TextValue $ws.Cells.Item(6, 4) "6.468"
Set-TextValue $ws.Cells.Item(7, 4) "0.8050"
Set-TextValue $ws.Cells.Item(8, 4) "1.058"
Set-TextValue $ws.Cells.Item(13, 4) "0.02931"
Set-TextValue $ws.Cells.Item(14, 4) "0.09261"
Set-TextValue $ws.Cells.Item(15, 4) "0.001676"
Set-TextValue $ws.Cells.Item(16, 4) "3.216"
Set-TextValue $ws.Cells.Item(17, 4) "0.04716"
Set-TextValue $ws.Cells.Item(27, 4) "0.0003304"
Set-TextValue $ws.Cells.Item(40, 4) "0.04182"
Set-TextValue $ws.Cells.Item(41, 4) "0.006847"
Set-TextValue $ws.Cells.Item(43, 4) "0.1040"
Set-TextValue $ws.Cells.Item(44, 4) "0.009855"
Set-TextValue $ws.Cells.Item(45, 4) "0.00005651"
Set-TextValue $ws.Cells.Item(47, 4) "0.6808"

# --- Rows whose ranking shifted: Coin (B), Link (C), Price (D), Volume(1h) (E) ---
$rows = @(
    @{ Row = 9;  Coin = "WazirX";                             Link = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx";                                           Price = "0.1431";    Volume = "8WazirXWRX" },
    @{ Row = 10; Coin = "MandalaExchangeToken";                Link = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx";                        Price = "0.07311";   Volume = "9MandalaExchangeTokenMDX" },
    @{ Row = 11; Coin = "LiechtensteinCryptoassetsExchange";   Link = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx";                Price = "0.03204";   Volume = "10LiechtensteinCryptoassetsExchangeLCX" },
    @{ Row = 12; Coin = "ProBitToken";                         Link = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob";                                     Price = "0.1316";    Volume = "11ProBitTokenPROBBestin24h" },
    @{ Row = 18; Coin = "One";                                 Link = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one";                                          Price = "0.0005826"; Volume = "17OneONEWorstin24h" },
    @{ Row = 19; Coin = "TigerCash";                           Link = "https://coinranking.com/coin/6hIn06L2+tigercash-tch";                                        Price = "0.006274";  Volume = "18TigerCashTCH" },
    @{ Row = 20; Coin = "BitKan";                              Link = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan";                                      Price = "0.001052";  Volume = "19BitKanKAN" },
    @{ Row = 21; Coin = "HotbitToken";                         Link = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb";                                  Price = "0.004120";  Volume = "20HotbitTokenHTB" },
    @{ Row = 22; Coin = "NitroEx";                             Link = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx";                                       Price = "0.0001501"; Volume = "21NitroExNTX" },
    @{ Row = 23; Coin = "LEO";                                 Link = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";                                          Price = "3.968";     Volume = "22LEOLEO" },
    @{ Row = 24; Coin = "GateToken";                           Link = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt";                                     Price = "3.378";     Volume = "23GateTokenGT" },
    @{ Row = 25; Coin = "BTSEToken";                           Link = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";                                   Price = "2.124";     Volume = "24BTSETokenBTSE" },
    @{ Row = 26; Coin = "BitpandaEcosystemToken";              Link = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best";                          Price = "0.3272";    Volume = "25BitpandaEcosystemTokenBEST" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Coin
    $ws.Cells.Item($r.Row, 3).Value = $r.Link
    Set-TextValue $ws.Cells.Item($r.Row, 4) $r.Price
    $ws.Cells.Item($r.Row, 5).Value = $r.Volume
}

# --- Row 48 (BOLO): Price + Volume(1h) label change ---
Set-TextValue $ws.Cells.Item(48, 4) "0.02537"
$ws.Cells.Item(48, 5).Value = "47BOLOBOLO"

Write-Output "Applied symbol list update."
